$p = $ppt.ActivePresentation

# --- 1. Update the cached "today" text of every datetimeFigureOut date
#        placeholder (slide master + every slide layout) from 7/29/2016
#        to 8/8/2016. We locate the placeholder by its PlaceholderFormat
#        type (ppPlaceholderDate = 16) rather than by name, because the
#        shape name differs across layouts.
$ppPlaceholderDate = 16
$newDate = "8/8/2016"

$m = $p.SlideMaster

for ($i = 1; $i -le $m.Shapes.Count; $i++) {
    $sh = $m.Shapes.Item($i)
    if ($sh.Type -eq 14 -and $sh.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
        $sh.TextFrame.TextRange.Text = $newDate
    }
}

for ($li = 1; $li -le $m.CustomLayouts.Count; $li++) {
    $cl = $m.CustomLayouts.Item($li)
    for ($i = 1; $i -le $cl.Shapes.Count; $i++) {
        $sh = $cl.Shapes.Item($i)
        if ($sh.Type -eq 14 -and $sh.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

# --- 2. Remove the stray "TextBox 17" (text "0.75") from slide 1 - the
#        proposed-model callout that the commit separates from the
#        implementation drawing.
$s = $p.Slides.Item(1)
$s.Shapes.Item("TextBox 17").Delete()
